$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column M ("sampling_date") holds text values that look like Excel serial
# date numbers (e.g. "14169"). Revise them to the actual ISO date strings,
# keeping the cells as plain text (not auto-converted Excel date values).
$dates = @{
    2  = "2008-10-17"
    3  = "2008-10-17"
    4  = "2008-10-17"
    5  = "2008-10-17"
    6  = "2008-10-17"
    7  = "2008-10-17"
    8  = "2008-10-17"
    9  = "2008-10-17"
    10 = "2008-10-17"
    11 = "2008-10-17"
    12 = "2008-10-17"
    13 = "2008-10-17"
    14 = "2008-04-15"
    15 = "2008-10-17"
    16 = "2008-10-17"
    17 = "2008-10-17"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Range("M$row")
    # Force text formatting first so the date-like string isn't reinterpreted
    # as a date serial number by Excel's automatic type detection.
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$row]
    # Clean up back to the default (unstyled) cell formatting, matching the
    # original workbook which had no explicit style on these data cells.
    $cell.ClearFormats()
}
